$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "RequestToChangeER": change the selected ER name and move selection
# ---------------------------------------------------------------------------
$wsReq = $wb.Worksheets.Item("RequestToChangeER")
$wsReq.Range("A2").Value = "Alyazia Khamis"
$wsReq.Range("B2").Select()

# ---------------------------------------------------------------------------
# Sheet "ERList": append a new ER (row 4) with login email + password
# ---------------------------------------------------------------------------
$wsList = $wb.Worksheets.Item("ERList")
$wsList.Range("A4").Value = "Abdulla Khalid"
$wsList.Range("B4").Value = "ishikite@domy.me"
$wsList.Range("C4").Value = "Test@123"

$wsList.Hyperlinks.Add($wsList.Range("B4"), "mailto:ishikite@domy.me")
$wsList.Hyperlinks.Add($wsList.Range("C4"), "mailto:Test@123")

# Hyperlinks.Add() stamps its own (redundant) style xf; re-point the cells at
# the very same "Hyperlink" style already used by B2/C2/B3/C3 so no spurious
# style record sticks around.
$wsList.Range("B4").Style = $wsList.Range("B2").Style
$wsList.Range("C4").Style = $wsList.Range("C2").Style

# ---------------------------------------------------------------------------
# Sheet "Sheet2": populate it with the same ER roster shown on ERList (minus
# the newly added 4th row) plus matching hyperlinks / column widths.
# ---------------------------------------------------------------------------
$wsSheet2 = $wb.Worksheets.Item("Sheet2")
$wsList.Range("A1:C3").Copy($wsSheet2.Range("A1"))

$wsSheet2.Hyperlinks.Add($wsSheet2.Range("B2"), "mailto:ertesting2he@gmail.com")
$wsSheet2.Hyperlinks.Add($wsSheet2.Range("C2"), "mailto:Test@123")
$wsSheet2.Hyperlinks.Add($wsSheet2.Range("B3"), "mailto:er3hadeel@gmail.com")
$wsSheet2.Hyperlinks.Add($wsSheet2.Range("C3"), "mailto:Test@123")

$wsSheet2.Range("B2").Style = $wsList.Range("B2").Style
$wsSheet2.Range("C2").Style = $wsList.Range("C2").Style
$wsSheet2.Range("B3").Style = $wsList.Range("B3").Style
$wsSheet2.Range("C3").Style = $wsList.Range("C3").Style

$wsSheet2.Columns.Item(1).ColumnWidth = 25
$wsSheet2.Columns.Item(2).ColumnWidth = 29.5

$wsSheet2.Range("A1:C4").Select()

# ---------------------------------------------------------------------------
# Make "ERList" the active tab/sheet (select its last populated cell) - do
# this last so it "wins" the single tabSelected flag across the workbook.
# ---------------------------------------------------------------------------
$wsList.Range("C4").Select()
$wsList.Activate()
